$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 21330051920031
$ws.Range("B2").Value = "AGUILAR"
$ws.Range("C2").Value = "ARIAS"
$ws.Range("D2").Value = "ALESSANDRA"
$ws.Range("E2").Value = "TECNOLOGÍAS DE LA INFORMACIÓN Y LA COMUNICACIÓN"
$ws.Range("F2").Value = "1BV"
$ws.Range("G2").Value = 6
